$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# Rushing sheet updates
$rushing.Range("C2").Value = 8
$rushing.Range("E2").Value = 12

$rushing.Range("C6").Value = 3

$rushing.Range("E7").Value = 3
$rushing.Range("F7").Value = 3

$rushing.Range("C8").Value = 39
$rushing.Range("D8").Value = 21
$rushing.Range("E8").Value = 5
$rushing.Range("F8").Value = 6

$rushing.Range("C9").Value = 18
$rushing.Range("D9").Value = 17
$rushing.Range("E9").Value = 5
$rushing.Range("F9").Value = 10

# Receiving sheet updates
$receiving.Range("C4").Value = 13
$receiving.Range("D4").Value = 11
$receiving.Range("G4").Value = 3
$receiving.Range("H4").Value = 2

$receiving.Range("C5").Value = 23
$receiving.Range("D5").Value = 21

$receiving.Range("C6").Value = 5
$receiving.Range("D6").Value = 4
$receiving.Range("G6").Value = 2
$receiving.Range("H6").Value = 1

$receiving.Range("C7").Value = 73
$receiving.Range("D7").Value = 42
$receiving.Range("E7").Value = 25
$receiving.Range("F7").Value = 18
$receiving.Range("G7").Value = 9

$receiving.Range("E8").Value = 19
$receiving.Range("F8").Value = 11

$receiving.Range("C10").Value = 15
$receiving.Range("G10").Value = 4

$receiving.Range("C11").Value = 21
$receiving.Range("D11").Value = 12
$receiving.Range("G11").Value = 5
$receiving.Range("H11").Value = 4

$receiving.Range("C12").Value = 80
$receiving.Range("D12").Value = 66
$receiving.Range("E12").Value = 26
$receiving.Range("F12").Value = 19
$receiving.Range("G12").Value = 6
$receiving.Range("H12").Value = 6

$wb.Save()
